# Weekly CompStat data refresh: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings) - edit only the characters
# that actually changed so the surrounding run formatting is left intact.
# ---------------------------------------------------------------------------

# A8: "Volume 32   Number  3" -> "Volume 32   Number  4"
$ws.Range("A8").Characters(21, 1).Text = "4"

# C9: "Report Covering the Week  1/13/2025  Through  1/19/2025"
#  -> "Report Covering the Week  1/20/2025  Through  1/26/2025"
$ws.Range("C9").Characters(27, 9).Text = "1/20/2025"
$ws.Range("C9").Characters(47, 9).Text = "1/26/2025"

# ---------------------------------------------------------------------------
# Helper: set a cell's value while preserving its current (numeric) style,
# used for the plain number -> number edits below.
# ---------------------------------------------------------------------------

# Row 14 - Murder
$ws.Range("C14").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("N14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -88.888888888888

# Row 15 - Rape
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 10
$ws.Range("H15").Value = -23.076923076923
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -10
$ws.Range("L15").Value = 12.5
$ws.Range("M15").Value = -10
$ws.Range("N15").Value = -52.631578947368

# Row 16 - Robbery
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 23
$ws.Range("E16").Value = 43.478260869565
$ws.Range("F16").Value = 110
$ws.Range("G16").Value = 117
$ws.Range("H16").Value = -5.982905982905
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = -23.846153846153
$ws.Range("M16").Value = -4.807692307692
$ws.Range("N16").Value = -86.835106382978

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 26
$ws.Range("D17").Value = 25
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 131
$ws.Range("G17").Value = 142
$ws.Range("H17").Value = -7.746478873239
$ws.Range("I17").Value = 116
$ws.Range("J17").Value = 128
$ws.Range("K17").Value = -9.375
$ws.Range("L17").Value = -21.088435374149
$ws.Range("M17").Value = 27.472527472527
$ws.Range("N17").Value = -51.464435146443

# Row 18 - Burglary
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 44
$ws.Range("E18").Value = -22.727272727272
$ws.Range("F18").Value = 127
$ws.Range("G18").Value = 160
$ws.Range("H18").Value = -20.625
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 152
$ws.Range("K18").Value = -23.684210526315
$ws.Range("L18").Value = -34.831460674157
$ws.Range("M18").Value = -36.263736263736
$ws.Range("N18").Value = -88.175331294597

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 156
$ws.Range("D19").Value = 232
$ws.Range("E19").Value = -32.758620689655
$ws.Range("F19").Value = 763
$ws.Range("G19").Value = 928
$ws.Range("H19").Value = -17.780172413793
$ws.Range("I19").Value = 694
$ws.Range("J19").Value = 823
$ws.Range("K19").Value = -15.674362089914
$ws.Range("L19").Value = -13.574097135741
$ws.Range("M19").Value = -8.443271767810
$ws.Range("N19").Value = -70.754319426885

# Row 20 - G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 33
$ws.Range("H20").Value = -39.393939393939
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -34.482758620689
$ws.Range("L20").Value = -38.709677419354
$ws.Range("M20").Value = 26.666666666666
$ws.Range("N20").Value = -96.049896049896

# Row 21 - TOTAL
$ws.Range("C21").Value = 254
$ws.Range("D21").Value = 330
$ws.Range("E21").Value = -23.030303030303
$ws.Range("F21").Value = 1162
$ws.Range("G21").Value = 1394
$ws.Range("H21").Value = -16.642754662840
$ws.Range("I21").Value = 1054
$ws.Range("J21").Value = 1253
$ws.Range("K21").Value = -15.881883479648
$ws.Range("L21").Value = -18.735543562066
$ws.Range("M21").Value = -9.294320137693
$ws.Range("N21").Value = -78.285949732179

# Row 22 - Transit
$ws.Range("C22").Value = 11
$ws.Range("D22").Value = 17
$ws.Range("E22").Value = -35.294117647058
$ws.Range("F22").Value = 45
$ws.Range("G22").Value = 64
$ws.Range("H22").Value = -29.6875
$ws.Range("I22").Value = 40
$ws.Range("J22").Value = 56
$ws.Range("K22").Value = -28.571428571428
$ws.Range("L22").Value = 2.564102564102
$ws.Range("M22").Value = -23.076923076923

# Row 23 - Housing
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 250
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = 71.428571428571
$ws.Range("I23").Value = 31
$ws.Range("J23").Value = 18
$ws.Range("K23").Value = 72.222222222222
$ws.Range("L23").Value = 34.782608695652
$ws.Range("M23").Value = 19.230769230769

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 414
$ws.Range("D24").Value = 402
$ws.Range("E24").Value = 2.985074626865
$ws.Range("F24").Value = 1593
$ws.Range("G24").Value = 1581
$ws.Range("H24").Value = 0.759013282732
$ws.Range("I24").Value = 1444
$ws.Range("J24").Value = 1435
$ws.Range("K24").Value = 0.627177700348
$ws.Range("L24").Value = 6.568265682656
$ws.Range("M24").Value = 26.223776223776

# Row 25 - Retail Theft
$ws.Range("C25").Value = 369
$ws.Range("D25").Value = 338
$ws.Range("E25").Value = 9.171597633136
$ws.Range("F25").Value = 1287
$ws.Range("G25").Value = 1338
$ws.Range("H25").Value = -3.811659192825
$ws.Range("I25").Value = 1168
$ws.Range("J25").Value = 1218
$ws.Range("K25").Value = -4.105090311986
$ws.Range("L25").Value = 1.919720767888

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 83
$ws.Range("D26").Value = 69
$ws.Range("E26").Value = 20.289855072463
$ws.Range("F26").Value = 352
$ws.Range("G26").Value = 352
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 314
$ws.Range("J26").Value = 318
$ws.Range("K26").Value = -1.257861635220
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 33.050847457627

# Row 27 - UCR Rape*
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 20
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 11
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = -31.25
$ws.Range("L27").Value = -26.666666666666

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -41.666666666666
$ws.Range("F28").Value = 53
$ws.Range("G28").Value = 50
$ws.Range("H28").Value = 6
$ws.Range("I28").Value = 46
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = 6.976744186046
$ws.Range("L28").Value = -13.207547169811

# Row 29 - Shooting Vic.
$ws.Range("F29").Value = "'0"
$ws.Range("I29").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100

# Row 30 - Shooting Inc.
$ws.Range("F30").Value = "'0"
$ws.Range("I30").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100

# Row 31 - Hate Crimes
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 33.333333333333
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = 33.333333333333
$ws.Range("L31").Value = -33.333333333333
